$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$newTimes = @(
    "2021-10-05 14:35:16.875245",
    "2021-10-05 14:35:16.875253",
    "2021-10-05 14:35:16.875256",
    "2021-10-05 14:35:16.875259",
    "2021-10-05 14:35:16.875262",
    "2021-10-05 14:35:16.875264",
    "2021-10-05 14:35:16.875267",
    "2021-10-05 14:35:16.875269",
    "2021-10-05 14:35:16.875272",
    "2021-10-05 14:35:16.875275",
    "2021-10-05 14:35:16.875277",
    "2021-10-05 14:35:16.875280",
    "2021-10-05 14:35:16.875282",
    "2021-10-05 14:35:16.875285",
    "2021-10-05 14:35:16.875287",
    "2021-10-05 14:35:16.875290",
    "2021-10-05 14:35:16.875293",
    "2021-10-05 14:35:16.875296",
    "2021-10-05 14:35:16.875298",
    "2021-10-05 14:35:16.875301",
    "2021-10-05 14:35:16.875303",
    "2021-10-05 14:35:16.875306",
    "2021-10-05 14:35:16.875308",
    "2021-10-05 14:35:16.875311",
    "2021-10-05 14:35:16.875314",
    "2021-10-05 14:35:16.875316",
    "2021-10-05 14:35:16.875319",
    "2021-10-05 14:35:16.875321",
    "2021-10-05 14:35:16.875324",
    "2021-10-05 14:35:16.875327",
    "2021-10-05 14:35:16.875329",
    "2021-10-05 14:35:16.875332",
    "2021-10-05 14:35:16.875335",
    "2021-10-05 14:35:16.875337",
    "2021-10-05 14:35:16.875340",
    "2021-10-05 14:35:16.875342",
    "2021-10-05 14:35:16.875345",
    "2021-10-05 14:35:16.875347",
    "2021-10-05 14:35:16.875350",
    "2021-10-05 14:35:16.875352",
    "2021-10-05 14:35:16.875355",
    "2021-10-05 14:35:16.875358",
    "2021-10-05 14:35:16.875360",
    "2021-10-05 14:35:16.875363",
    "2021-10-05 14:35:16.875365",
    "2021-10-05 14:35:16.875368",
    "2021-10-05 14:35:16.875370",
    "2021-10-05 14:35:16.875373",
    "2021-10-05 14:35:16.875375",
    "2021-10-05 14:35:16.875378",
    "2021-10-05 14:35:16.875381",
    "2021-10-05 14:35:16.875383",
    "2021-10-05 14:35:16.875386",
    "2021-10-05 14:35:16.875389",
    "2021-10-05 14:35:16.875426",
    "2021-10-05 14:35:16.875440",
    "2021-10-05 14:35:16.875446",
    "2021-10-05 14:35:16.875451",
    "2021-10-05 14:35:16.875455",
    "2021-10-05 14:35:16.875460",
    "2021-10-05 14:35:16.875465",
    "2021-10-05 14:35:16.875468",
    "2021-10-05 14:35:16.875470",
    "2021-10-05 14:35:16.875473",
    "2021-10-05 14:35:16.875478",
    "2021-10-05 14:35:16.875481",
    "2021-10-05 14:35:16.875484",
    "2021-10-05 14:35:16.875486",
    "2021-10-05 14:35:16.875489",
    "2021-10-05 14:35:16.875491",
    "2021-10-05 14:35:16.875494",
    "2021-10-05 14:35:16.875497",
    "2021-10-05 14:35:16.875499",
    "2021-10-05 14:35:16.875502",
    "2021-10-05 14:35:16.875505"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add the new "metadata" sheet right after "data" ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# D2 must hold the literal text "0.164" (not the number 0.164). Force text
# via a temporary Text number format, then clear the format again so the
# cell ends up with no explicit style, matching a freshly authored sheet.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.164"
$ws.Range("D2").ClearFormats()

# Copy the header style (bold, bordered, centered) from the data sheet header row
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Polymicrogyria and Schizencephaly"
$ws.Range("C2").Value = 18
$ws.Range("E2").Value = "2021-07-08T03:47:57.928004Z"
$ws.Range("F2").Value = "2021-10-05 14:35:16.871915"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/18/?format=json"

